# Add a new statistic/record row (row 3) to the "Verkenners" sheet,
# matching the layout: ID | Kind Naam en Van | Ouer Naam en Van | Selfoon Nommer | Diens
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "TestName TestSurname"
$ws.Range("C3").Value = "TestParentName TestParent Surname"
$ws.Range("D3").Value = "085 555 5555"
$ws.Range("E3").Value = "Eerste"
